# feat: start page ftp
#
# Adds a new "ftp" worksheet between "entidade" and "setor" describing the
# VTEX FTP import/export folder layout, and makes a couple of small edits
# on the "entidade" sheet (a running counter bump + a code/CNPJ fix).

$wb = $excel.ActiveWorkbook
$entidade = $wb.Worksheets.Item("entidade")
$setor = $wb.Worksheets.Item("setor")

# ---------------------------------------------------------------------
# 1. entidade sheet tweaks
# ---------------------------------------------------------------------
$entidade.Range("A2").Value = "x1"
$entidade.Range("B2").Value = "x1"
$entidade.Range("C2").Value = "07.046.881/1007-00"
$entidade.Range("E2").Value = 144

# ---------------------------------------------------------------------
# 2. Insert the new "ftp" sheet right after "entidade"
# ---------------------------------------------------------------------
$ftp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $entidade)
$ftp.Name = "ftp"

# Column order on the sheet is A,B,C,D,E (importacao / backup_importacao /
# exportacao / backup_exportacao / erro) but the values are entered in the
# order below so new entries land in the shared-string table the same way
# they did originally.
$ftp.Range("A1").Value = "importacao"
$ftp.Range("C1").Value = "exportacao"
$ftp.Range("B1").Value = "backup_importacao"
$ftp.Range("D1").Value = "backup_exportacao"
$ftp.Range("E1").Value = "erro"

$ftp.Range("A2").Value = "/home/ftpsynapcomp/Embu/Vtex/importacao"
$ftp.Range("B2").Value = "/home/ftpsynapcomp/Embu/Vtex/bkp_importacao"
$ftp.Range("C2").Value = "/home/ftpsynapcomp/Embu/Vtex/exportacao"
$ftp.Range("D2").Value = "/home/ftpsynapcomp/Embu/Vtex/bkp_exportacao"
$ftp.Range("E2").Value = "/home/ftpsynapcomp/Embu/Vtex/erro"

$ftp.Columns.Item(1).ColumnWidth = 42.166666666666664
$ftp.Columns.Item(2).ColumnWidth = 45.736979166666664
$ftp.Columns.Item(3).ColumnWidth = 40.307291666666664
$ftp.Columns.Item(4).ColumnWidth = 46.022135416666664
$ftp.Columns.Item(5).ColumnWidth = 45.451822916666664

# ---------------------------------------------------------------------
# 3. Selections: entidade -> G6 (single cell), ftp -> I8 (ends active)
#    "setor" selection/view is left exactly as it was.
# ---------------------------------------------------------------------
$entidade.Activate()
$entidade.Range("G6").Select()

$ftp.Activate()
$ftp.Range("I8").Select()
